$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns("N").Insert()

# Match the new column's width to the neighbouring "In Advance" column (M)
$wsSchedule.Columns("N").ColumnWidth = $wsSchedule.Columns("M").ColumnWidth

# --- Make "Repayment schedule" the active sheet/tab ---
$wsSchedule.Activate()
$wsSchedule.Range("K14").Select()
